$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): rename D1/E1/F1 and add new choice columns G1..N1 ---
$ws.Range("D1").Value = "slider_response"
$ws.Range("E1").Value = "choice_1"
$ws.Range("F1").Value = "choice_2"

# New header cells G1..N1 need the same style as the existing header cells (s="1").
# Copy style from F1 (already s="1") into each new cell before setting its value,
# so we reuse the existing cellXf instead of Excel fabricating a brand-new one.
$ws.Range("F1").Copy($ws.Range("G1"))
$ws.Range("G1").Value = "choice_3"

$ws.Range("F1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "choice_4"

$ws.Range("F1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "choice_5"

$ws.Range("F1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "choice_6"

$ws.Range("F1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "choice_7"

$ws.Range("F1").Copy($ws.Range("L1"))
$ws.Range("L1").Value = "choice_8"

$ws.Range("F1").Copy($ws.Range("M1"))
$ws.Range("M1").Value = "choice_9"

$ws.Range("F1").Copy($ws.Range("N1"))
$ws.Range("N1").Value = "choice_10"

# --- Row 2: add empty styled cell H2 (style copied from the existing G2) ---
$ws.Range("G2").Copy($ws.Range("H2"))

# --- Row 3: add empty styled cell E3 (style copied from the existing D3) ---
$ws.Range("D3").Copy($ws.Range("E3"))

# --- Column widths: column F joins the 18.29-wide D:E block ---
$ws.Range("F1").ColumnWidth = 17.5
